# Nieuwe data toegevoegd via Streamlit op 2024-12-03 18:10:14
# Append one new inspection record as row 75 on the (single) data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 75

$ws.Range("A$newRow").Value = "Kindergarden"
$ws.Range("B$newRow").Value = "Kindergarden Hilversum Johannes Geradtsweg"
$ws.Range("C$newRow").Value = "KDV"

# Force text storage so the report date stays a literal "YYYY-MM-DD" string
# (matching the rest of the sheet) instead of being auto-converted to a
# serial date value.
$ws.Range("D$newRow").NumberFormat = "@"
$ws.Range("D$newRow").Value = "2024-03-20"

$ws.Range("E$newRow").Value = 0
$ws.Range("F$newRow").Value = 0
$ws.Range("G$newRow").Value = 0
$ws.Range("H$newRow").Value = 1
$ws.Range("I$newRow").Value = 0
$ws.Range("J$newRow").Value = 0
